$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 990111.3636582257
$ws.Range("C9").Value = 831314.4975979078
$ws.Range("D9").Value = 1821425.861256134
$ws.Range("B10").Value = 1009838.458882504
$ws.Range("C10").Value = 850373.6581766765
$ws.Range("D10").Value = 1860212.117059181
$ws.Range("B11").Value = 1027531.694762051
$ws.Range("C11").Value = 868156.3520049311
$ws.Range("D11").Value = 1895688.046766982
$ws.Range("B12").Value = 1038369.02837113
$ws.Range("C12").Value = 878556.1803677904
$ws.Range("D12").Value = 1916925.20873892
$ws.Range("B13").Value = 1037644.072918339
$ws.Range("C13").Value = 880302.1978995795
$ws.Range("D13").Value = 1917946.270817919
$ws.Range("B14").Value = 1031218.814306257
$ws.Range("C14").Value = 877441.7585444573
$ws.Range("D14").Value = 1908660.572850714
$ws.Range("B15").Value = 1018946.673769818
$ws.Range("C15").Value = 866333.4608667413
$ws.Range("D15").Value = 1885280.134636559
$ws.Range("B16").Value = 999525.3877778889
$ws.Range("C16").Value = 846011.7320133122
$ws.Range("D16").Value = 1845537.119791201
$ws.Range("B17").Value = 972793.8281127313
$ws.Range("C17").Value = 821655.2534553671
$ws.Range("D17").Value = 1794449.081568098
$ws.Range("B18").Value = 943897.9557722975
$ws.Range("C18").Value = 795191.0654827683
$ws.Range("D18").Value = 1739089.021255066
$ws.Range("B19").Value = 909182.543202111
$ws.Range("C19").Value = 766968.9717814658
$ws.Range("D19").Value = 1676151.514983577
$ws.Range("B20").Value = 864782.6833214258
$ws.Range("C20").Value = 728187.7033590564
$ws.Range("D20").Value = 1592970.386680482
$ws.Range("B21").Value = 814816.909249138
$ws.Range("C21").Value = 686234.4709944156
$ws.Range("D21").Value = 1501051.380243554
$ws.Range("B22").Value = 763954.3121890175
$ws.Range("C22").Value = 643899.7470905373
$ws.Range("D22").Value = 1407854.059279555
$ws.Range("B23").Value = 708329.7785302646
$ws.Range("C23").Value = 595847.0536259327
$ws.Range("D23").Value = 1304176.832156197
$ws.Range("B24").Value = 648565.7703387433
$ws.Range("C24").Value = 545519.484876293
$ws.Range("D24").Value = 1194085.255215036
